# posts.xlsx was edited to drop the "「不可能はない」" post entry (row 823).
# Removing the entire row shifts every subsequent row (824-838) up by one,
# which also shrinks the sheet's used range from A1:C838 to A1:C837.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("823").Delete()
